$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create row 57 with the same formatting as row 56 (copy style of A56 -> A57)
$ws.Range("A56").Copy()
$ws.Range("A57").PasteSpecial(-4122)

$data = @(
    @(0, 33.94444444444444, 1.95, 0, 0.138),
    @(1, 33.94444444444444, 1.95, 2, 0.001),
    @(2, 33.94444444444444, 1.95, 3, 0.003),
    @(3, 33.94444444444444, 1.95, 4, 0.009000000000000001),
    @(4, 33.94444444444444, 1.95, 5, 0.023),
    @(5, 33.94444444444444, 1.95, 6, 0.041),
    @(6, 33.94444444444444, 1.95, 7, 0.058),
    @(7, 33.94444444444444, 1.95, 8, 0.039),
    @(8, 33.94444444444444, 1.95, 9, 0.038),
    @(9, 33.94444444444444, 1.95, 10, 0.038),
    @(10, 33.94444444444444, 1.95, 11, 0.033),
    @(11, 33.94444444444444, 1.95, 12, 0.02),
    @(12, 33.94444444444444, 1.95, 13, 0.024),
    @(13, 33.94444444444444, 1.95, 14, 0.036),
    @(14, 33.94444444444444, 1.95, 15, 0.037),
    @(15, 33.94444444444444, 1.95, 16, 0.042),
    @(16, 33.94444444444444, 1.95, 17, 0.033),
    @(17, 33.94444444444444, 1.95, 18, 0.028),
    @(18, 33.94444444444444, 1.95, 19, 0.034),
    @(19, 33.94444444444444, 1.95, 20, 0.024),
    @(20, 33.94444444444444, 1.95, 21, 0.022),
    @(21, 33.94444444444444, 1.95, 22, 0.016),
    @(22, 33.94444444444444, 1.95, 23, 0.018),
    @(23, 33.94444444444444, 1.95, 24, 0.032),
    @(24, 33.94444444444444, 1.95, 25, 0.022),
    @(25, 33.94444444444444, 1.95, 26, 0.027),
    @(26, 33.94444444444444, 1.95, 27, 0.026),
    @(27, 33.94444444444444, 1.95, 28, 0.009000000000000001),
    @(28, 33.94444444444444, 1.95, 29, 0.018),
    @(29, 33.94444444444444, 1.95, 30, 0.009000000000000001),
    @(30, 33.94444444444444, 1.95, 31, 0.011),
    @(31, 33.94444444444444, 1.95, 32, 0.011),
    @(32, 33.94444444444444, 1.95, 33, 0.014),
    @(33, 33.94444444444444, 1.95, 34, 0.009000000000000001),
    @(34, 33.94444444444444, 1.95, 35, 0.007),
    @(35, 33.94444444444444, 1.95, 36, 0.003),
    @(36, 33.94444444444444, 1.95, 37, 0.005),
    @(37, 33.94444444444444, 1.95, 38, 0.002),
    @(38, 33.94444444444444, 1.95, 39, 0.005),
    @(39, 33.94444444444444, 1.95, 40, 0.008),
    @(40, 33.94444444444444, 1.95, 41, 0.003),
    @(41, 33.94444444444444, 1.95, 42, 0.003),
    @(42, 33.94444444444444, 1.95, 43, 0.002),
    @(43, 33.94444444444444, 1.95, 45, 0.001),
    @(44, 33.94444444444444, 1.95, 46, 0.001),
    @(45, 33.94444444444444, 1.95, 47, 0.002),
    @(46, 33.94444444444444, 1.95, 48, 0.002),
    @(47, 33.94444444444444, 1.95, 49, 0.002),
    @(48, 33.94444444444444, 1.95, 50, 0.001),
    @(49, 33.94444444444444, 1.95, 51, 0.001),
    @(50, 33.94444444444444, 1.95, 52, 0.003),
    @(51, 33.94444444444444, 1.95, 54, 0.001),
    @(52, 33.94444444444444, 1.95, 55, 0.001),
    @(53, 33.94444444444444, 1.95, 61, 0.001),
    @(54, 33.94444444444444, 1.95, 63, 0.001),
    @(55, 33.94444444444444, 1.95, 74, 0.001)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}
